$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 2).Value = 5.45195
$ws.Cells.Item(2, 3).Value = 5.89945
$ws.Cells.Item(2, 4).Value = 37.50169999999999
$ws.Cells.Item(2, 5).Value = 0.17115
$ws.Cells.Item(2, 6).Value = 2219.401700000001
$ws.Cells.Item(2, 7).Value = 2014.1057
$ws.Cells.Item(2, 8).Value = 205.29605
$ws.Cells.Item(2, 9).Value = 205.29605
$ws.Cells.Item(2, 11).Value = 2798.105
$ws.Cells.Item(2, 12).Value = 2014.125
$ws.Cells.Item(2, 13).Value = 783.98
$ws.Cells.Item(2, 14).Value = 345.6012
$ws.Cells.Item(2, 15).Value = 438.3792
$ws.Cells.Item(3, 2).Value = 5.789
$ws.Cells.Item(3, 3).Value = 9.571999999999999
$ws.Cells.Item(3, 4).Value = 28.079
$ws.Cells.Item(3, 6).Value = 2327.595
$ws.Cells.Item(3, 7).Value = 2122.899
$ws.Cells.Item(3, 8).Value = 204.696
$ws.Cells.Item(3, 9).Value = 193.677
$ws.Cells.Item(3, 10).Value = 11.019
$ws.Cells.Item(3, 11).Value = 2476.9342
$ws.Cells.Item(3, 12).Value = 2122.95
$ws.Cells.Item(3, 13).Value = 353.9842
$ws.Cells.Item(3, 14).Value = 308.7674
$ws.Cells.Item(3, 15).Value = 45.2168
$ws.Cells.Item(4, 2).Value = 6.746
$ws.Cells.Item(4, 3).Value = 8.249000000000001
$ws.Cells.Item(4, 4).Value = 41.042
$ws.Cells.Item(4, 6).Value = 2510.203
$ws.Cells.Item(4, 7).Value = 2280.003
$ws.Cells.Item(4, 8).Value = 230.2
$ws.Cells.Item(4, 9).Value = 64.40900000000001
$ws.Cells.Item(4, 10).Value = 165.791
$ws.Cells.Item(4, 11).Value = 2467.796
$ws.Cells.Item(4, 12).Value = 2280.003
$ws.Cells.Item(4, 13).Value = 187.793
$ws.Cells.Item(4, 14).Value = 176.6258
$ws.Cells.Item(4, 15).Value = 11.1674

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 2).Value = 5.591900000000001
$ws.Cells.Item(2, 3).Value = 6.396350000000001
$ws.Cells.Item(2, 4).Value = 36.86314999999999
$ws.Cells.Item(2, 5).Value = 0.1292
$ws.Cells.Item(2, 6).Value = 2253.7439
$ws.Cells.Item(2, 7).Value = 2044.3985
$ws.Cells.Item(2, 8).Value = 209.34555
$ws.Cells.Item(2, 9).Value = 209.34555
$ws.Cells.Item(2, 11).Value = 2285.459
$ws.Cells.Item(2, 12).Value = 2044.411
$ws.Cells.Item(2, 13).Value = 241.048
$ws.Cells.Item(2, 14).Value = 241.048
$ws.Cells.Item(3, 2).Value = 5.963
$ws.Cells.Item(3, 3).Value = 9.461
$ws.Cells.Item(3, 4).Value = 28.427
$ws.Cells.Item(3, 6).Value = 2359.512
$ws.Cells.Item(3, 7).Value = 2141.396
$ws.Cells.Item(3, 8).Value = 218.116
$ws.Cells.Item(3, 9).Value = 201.182
$ws.Cells.Item(3, 10).Value = 16.934
$ws.Cells.Item(3, 11).Value = 2347.4208
$ws.Cells.Item(3, 12).Value = 2141.429
$ws.Cells.Item(3, 13).Value = 205.9918
$ws.Cells.Item(3, 14).Value = 205.9918
$ws.Cells.Item(4, 2).Value = 6.7
$ws.Cells.Item(4, 3).Value = 8.093999999999999
$ws.Cells.Item(4, 4).Value = 47.574
$ws.Cells.Item(4, 6).Value = 2580.234
$ws.Cells.Item(4, 7).Value = 2317.715
$ws.Cells.Item(4, 8).Value = 262.519
$ws.Cells.Item(4, 9).Value = 54.052
$ws.Cells.Item(4, 10).Value = 208.467
$ws.Cells.Item(4, 11).Value = 2395.058
$ws.Cells.Item(4, 12).Value = 2317.715
$ws.Cells.Item(4, 13).Value = 77.343
$ws.Cells.Item(4, 14).Value = 77.343

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 2).Value = 5.43975
$ws.Cells.Item(2, 3).Value = 6.10355
$ws.Cells.Item(2, 4).Value = 37.0707
$ws.Cells.Item(2, 5).Value = 0.17115
$ws.Cells.Item(2, 6).Value = 2222.8903
$ws.Cells.Item(2, 7).Value = 2017.33845
$ws.Cells.Item(2, 8).Value = 205.55185
$ws.Cells.Item(2, 9).Value = 205.55185
$ws.Cells.Item(2, 11).Value = 2804.228
$ws.Cells.Item(2, 12).Value = 2017.357
$ws.Cells.Item(2, 13).Value = 786.8710000000001
$ws.Cells.Item(2, 14).Value = 359.913
$ws.Cells.Item(2, 15).Value = 426.9582
$ws.Cells.Item(3, 2).Value = 5.483
$ws.Cells.Item(3, 3).Value = 7.719
$ws.Cells.Item(3, 4).Value = 32.87
$ws.Cells.Item(3, 6).Value = 2320.004
$ws.Cells.Item(3, 7).Value = 2051.97
$ws.Cells.Item(3, 8).Value = 268.033
$ws.Cells.Item(3, 9).Value = 220.571
$ws.Cells.Item(3, 10).Value = 47.462
$ws.Cells.Item(3, 11).Value = 2618.863
$ws.Cells.Item(3, 12).Value = 2051.924
$ws.Cells.Item(3, 13).Value = 566.939
$ws.Cells.Item(3, 14).Value = 353.6808
$ws.Cells.Item(3, 15).Value = 213.2584
$ws.Cells.Item(4, 2).Value = 6.231
$ws.Cells.Item(4, 3).Value = 8.112
$ws.Cells.Item(4, 4).Value = 43.556
$ws.Cells.Item(4, 6).Value = 2573.291
$ws.Cells.Item(4, 7).Value = 2233.933
$ws.Cells.Item(4, 8).Value = 339.358
$ws.Cells.Item(4, 9).Value = 119.435
$ws.Cells.Item(4, 10).Value = 219.923
$ws.Cells.Item(4, 11).Value = 2490.067
$ws.Cells.Item(4, 12).Value = 2233.933
$ws.Cells.Item(4, 13).Value = 256.134
$ws.Cells.Item(4, 14).Value = 190.343
$ws.Cells.Item(4, 15).Value = 65.791

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 2).Value = 5.734050000000002
$ws.Cells.Item(2, 3).Value = 6.146849999999999
$ws.Cells.Item(2, 4).Value = 39.35619999999999
$ws.Cells.Item(2, 5).Value = 0.04355
$ws.Cells.Item(2, 6).Value = 2295.672100000001
$ws.Cells.Item(2, 7).Value = 2069.291249999999
$ws.Cells.Item(2, 8).Value = 226.3809
$ws.Cells.Item(2, 9).Value = 226.3809
$ws.Cells.Item(2, 11).Value = 2185.4816
$ws.Cells.Item(2, 12).Value = 2069.311
$ws.Cells.Item(2, 13).Value = 116.1706
$ws.Cells.Item(2, 14).Value = 80.70360000000001
$ws.Cells.Item(2, 15).Value = 35.467
$ws.Cells.Item(3, 2).Value = 6.113
$ws.Cells.Item(3, 3).Value = 9.366
$ws.Cells.Item(3, 4).Value = 29.606
$ws.Cells.Item(3, 6).Value = 2390.136
$ws.Cells.Item(3, 7).Value = 2164.023
$ws.Cells.Item(3, 8).Value = 226.112
$ws.Cells.Item(3, 9).Value = 213.845
$ws.Cells.Item(3, 10).Value = 12.267
$ws.Cells.Item(3, 11).Value = 2231.9788
$ws.Cells.Item(3, 12).Value = 2163.984
$ws.Cells.Item(3, 13).Value = 67.9948
$ws.Cells.Item(3, 14).Value = 67.9948
$ws.Cells.Item(4, 2).Value = 6.7
$ws.Cells.Item(4, 3).Value = 8.093999999999999
$ws.Cells.Item(4, 4).Value = 47.574
$ws.Cells.Item(4, 6).Value = 2580.234
$ws.Cells.Item(4, 7).Value = 2317.715
$ws.Cells.Item(4, 8).Value = 262.519
$ws.Cells.Item(4, 9).Value = 54.052
$ws.Cells.Item(4, 10).Value = 208.467
$ws.Cells.Item(4, 11).Value = 2328.9178
$ws.Cells.Item(4, 12).Value = 2317.715
$ws.Cells.Item(4, 13).Value = 11.2028
$ws.Cells.Item(4, 14).Value = 11.2028

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 2).Value = 5.656750000000001
$ws.Cells.Item(2, 3).Value = 5.934599999999999
$ws.Cells.Item(2, 4).Value = 39.77945
$ws.Cells.Item(2, 5).Value = 0.16955
$ws.Cells.Item(2, 6).Value = 2266.771
$ws.Cells.Item(2, 7).Value = 2056.2701
$ws.Cells.Item(2, 8).Value = 210.50085
$ws.Cells.Item(2, 9).Value = 210.50085
$ws.Cells.Item(2, 11).Value = 2393.0558
$ws.Cells.Item(2, 12).Value = 2056.285
$ws.Cells.Item(2, 13).Value = 336.7708
$ws.Cells.Item(2, 14).Value = 152.2866
$ws.Cells.Item(2, 15).Value = 184.4844
$ws.Cells.Item(3, 2).Value = 6.099
$ws.Cells.Item(3, 3).Value = 7.875
$ws.Cells.Item(3, 4).Value = 33.748
$ws.Cells.Item(3, 6).Value = 2371.432
$ws.Cells.Item(3, 7).Value = 2135.853
$ws.Cells.Item(3, 8).Value = 235.579
$ws.Cells.Item(3, 9).Value = 199.385
$ws.Cells.Item(3, 10).Value = 36.195
$ws.Cells.Item(3, 11).Value = 2332.4364
$ws.Cells.Item(3, 12).Value = 2135.884
$ws.Cells.Item(3, 13).Value = 196.5524
$ws.Cells.Item(3, 14).Value = 126.9988
$ws.Cells.Item(3, 15).Value = 69.5534
$ws.Cells.Item(4, 2).Value = 7.162
$ws.Cells.Item(4, 3).Value = 7.283
$ws.Cells.Item(4, 4).Value = 49.794
$ws.Cells.Item(4, 6).Value = 2576.937
$ws.Cells.Item(4, 7).Value = 2356.575
$ws.Cells.Item(4, 8).Value = 220.362
$ws.Cells.Item(4, 9).Value = 121.688
$ws.Cells.Item(4, 10).Value = 98.675
$ws.Cells.Item(4, 11).Value = 2467.0234
$ws.Cells.Item(4, 12).Value = 2356.575
$ws.Cells.Item(4, 13).Value = 110.4484
$ws.Cells.Item(4, 14).Value = 44.6264
$ws.Cells.Item(4, 15).Value = 65.8222
